# Applies the "Updated cryptos list" price/volume refresh (and the three
# Coin/Link row swaps) from the commit diff onto Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cellRef -> new text; cells whose text looks like a plain number are prefixed
# with a leading apostrophe so Excel's Value setter keeps them as text (matching
# the inlineStr cells in the workbook) instead of silently parsing them as numbers
# (which would drop things like trailing zeros, e.g. "1.10" -> 1.1).
$updates = [ordered]@{
    "D2" = "61.969.84"
    "E2" = "  -0.25%  "
    "D3" = "3.419.42"
    "E3" = "  -0.59%  "
    "E4" = "  -0.03%  "
    "D5" = "'412.07"
    "E5" = "  +0.63%  "
    "D6" = "'129.33"
    "E6" = "  +0.05%  "
    "D7" = "'0.630"
    "E7" = "  +0.77%  "
    "D9" = "'0.729"
    "E9" = "  -1.81%  "
    "E10" = "  -1.75%  "
    "D11" = "'43.23"
    "E11" = "  +0.75%  "
    "D12" = "'9.17"
    "E12" = "  +2.42%  "
    "D13" = "3.957.03"
    "E13" = "  -0.49%  "
    "E14" = "  +0.19%  "
    "E15" = "  +3.05%  "
    "D16" = "'20.99"
    "E16" = "  -1.67%  "
    "D17" = "3.422.19"
    "E17" = "  -1.69%  "
    "B18" = "Polygon"
    "C18" = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
    "D18" = "'1.10"
    "E18" = "  +2.85%  "
    "B19" = "Uniswap"
    "C19" = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
    "D19" = "'12.55"
    "E19" = "  +0.84%  "
    "D20" = "61.900.90"
    "E20" = "  -0.12%  "
    "D21" = "'488.36"
    "E21" = "  +20.76%  "
    "D22" = "'91.84"
    "E22" = "  +1.94%  "
    "D23" = "'3.30"
    "E23" = "  +2.98%  "
    "D24" = "'13.44"
    "E24" = "  -0.06%  "
    "D25" = "'3.36"
    "E25" = "  +4.31%  "
    "D26" = "'34.28"
    "E26" = "  +2.51%  "
    "D27" = "'9.19"
    "E27" = "  +4.71%  "
    "E28" = "  +0.27%  "
    "D29" = "'7.78"
    "E29" = "  +2.23%  "
    "B30" = "Cosmos"
    "C30" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D30" = "'12.03"
    "E30" = "  +0.47%  "
    "B31" = "Toncoin"
    "C31" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D31" = "'2.68"
    "E31" = "  -3.32%  "
    "E32" = "  -2.53%  "
    "E33" = "  -4.38%  "
    "D34" = "'41.88"
    "E34" = "  -4.10%  "
    "D35" = "'1.00"
    "E35" = "  +0.02%  "
    "D36" = "'58.31"
    "E36" = "  +7.27%  "
    "D37" = "'0.0493"
    "E37" = "  -2.34%  "
    "E38" = "  -0.06%  "
    "D39" = "'151.44"
    "E39" = "  +7.10%  "
    "D40" = "'3.42"
    "E40" = "  +0.66%  "
    "D41" = "'0.137"
    "E41" = "  +3.13%  "
    "B42" = "ARBITRUM"
    "C42" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D42" = "'2.13"
    "E42" = "  +6.90%  "
    "B43" = "TheGraph"
    "C43" = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
    "D43" = "'0.322"
    "E43" = "  +2.73%  "
    "E44" = "  +0.57%  "
    "D45" = "'2.64"
    "E45" = "  +10.28%  "
    "D46" = "'4.23"
    "E46" = "  +4.47%  "
    "D47" = "'2.35"
    "E47" = "  +21.21%  "
    "D48" = "'16.52"
    "E48" = "  -1.45%  "
    "B49" = "BitcoinSV"
    "C49" = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
    "D49" = "'117.70"
    "E49" = "  +21.73%  "
    "B50" = "EnergySwap"
    "C50" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D50" = "'22.62"
    "E50" = "  +2.44%  "
    "D51" = "'0.149"
    "E51" = "  +15.63%  "
}

# A cell known to carry the default (unstyled) format, used below to strip the
# quote-prefix style Excel applies when it text-quotes a numeric-looking value,
# so the edited cells keep the same (style-less) formatting as before the edit.
$defaultStyle = $ws.Range("B2").Style

foreach ($ref in $updates.Keys) {
    $text = $updates[$ref]
    $ws.Range($ref).Value = $text
    if ($text.StartsWith("'")) {
        $ws.Range($ref).Style = $defaultStyle
    }
}
